$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.107940196990967
$ws.Range("B1").Value = 1.266552925109863
$ws.Range("C1").Value = 1.590944051742554
$ws.Range("D1").Value = 3.254141569137573
$ws.Range("E1").Value = 4.092807769775391
